# Update gh-pages to output generated at 456a3b4
# Applies refreshed "F" column (number of comments/participants) values
# across the "展览" (sheet 1), "本地生活" (sheet 3) and "全部类型" (sheet 4)
# worksheets, leaving every other cell untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value  = 9448
$ws1.Range("F8").Value  = 7418
$ws1.Range("F12").Value = 32
$ws1.Range("F13").Value = 6725
$ws1.Range("F16").Value = 459
$ws1.Range("F18").Value = 662
$ws1.Range("F24").Value = 10928
$ws1.Range("F26").Value = 54
$ws1.Range("F27").Value = 2074
$ws1.Range("F28").Value = 2679
$ws1.Range("F31").Value = 2424
$ws1.Range("F34").Value = 34
$ws1.Range("F37").Value = 1506
$ws1.Range("F40").Value = 5540
$ws1.Range("F42").Value = 785
$ws1.Range("F43").Value = 143
$ws1.Range("F49").Value = 1115

# --- Sheet 3: 本地生活 -------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 226

# --- Sheet 4: 全部类型 -------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value  = 9448
$ws4.Range("F8").Value  = 226
$ws4.Range("F11").Value = 7421
$ws4.Range("F14").Value = 32
$ws4.Range("F15").Value = 6725
$ws4.Range("F16").Value = 6725
$ws4.Range("F19").Value = 459
$ws4.Range("F20").Value = 662
$ws4.Range("F26").Value = 205
$ws4.Range("F27").Value = 10928
$ws4.Range("F29").Value = 54
$ws4.Range("F30").Value = 2074
$ws4.Range("F31").Value = 2679
$ws4.Range("F32").Value = 2424
$ws4.Range("F36").Value = 34
$ws4.Range("F39").Value = 1506
$ws4.Range("F40").Value = 5540
$ws4.Range("F43").Value = 785
$ws4.Range("F44").Value = 143
$ws4.Range("F50").Value = 1115

Write-Host "Applied 38 cell updates across 3 worksheets"
